# Update "want to go" counts (column F) for several events on both the
# "展览" (Exhibition) sheet and the aggregated "全部类型" (All Types) sheet.
# These sheets contain duplicate rows for the same events, so the same
# numeric bumps are applied in both places.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Row -> new value map (column F) for the "展览" sheet
$exhibitUpdates = @{
    12 = 742
    14 = 1875
    15 = 393
    16 = 4667
    17 = 397
    21 = 157
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Range("F$row").Value = $exhibitUpdates[$row]
}

# Row -> new value map (column F) for the "全部类型" sheet
$allUpdates = @{
    26 = 742
    29 = 1875
    30 = 393
    31 = 4667
    33 = 397
    38 = 157
}

foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}
